# Append two more user-credential rows to the "Sheet" worksheet, matching
# the new entries added in the commit ("I have complete flipkart login program").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Raju123"
$ws.Range("B4").Value = "Raju@12345"

$ws.Range("A5").Value = "asd"
$ws.Range("B5").Value = "Asd@123"
